# Update cryptocurrency price/volume figures to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell, new text value) pairs taken from the updated symbol list.
# Each pair is prefixed with the unary "," so PowerShell keeps it as a
# nested 2-element array instead of flattening everything into one list.
$updates = @(
    ,@('D2', '261.53')
    ,@('E2', '0.24%')
    ,@('D3', '26.67')
    ,@('E3', '-2.38%')
    ,@('D4', '4.702')
    ,@('E4', '0.01%')
    ,@('D5', '0.06184')
    ,@('D6', '6.709')
    ,@('E6', '0.55%')
    ,@('D7', '0.8503')
    ,@('E7', '0.36%')
    ,@('D8', '0.9120')
    ,@('E8', '-2.26%')
    ,@('E9', '-0.25%')
    ,@('D10', '0.05134')
    ,@('E10', '4.79%')
    ,@('D11', '0.07098')
    ,@('E11', '-0.06%')
    ,@('D12', '0.03109')
    ,@('E12', '1.10%')
    ,@('D13', '0.09037')
    ,@('E13', '-0.30%')
    ,@('D14', '0.001544')
    ,@('E14', '0.37%')
    ,@('D15', '0.0006149')
    ,@('E15', '0.88%')
    ,@('D16', '0.005988')
    ,@('E16', '-1.34%')
    ,@('D17', '3.449')
    ,@('E17', '0.00%')
    ,@('D18', '3.170')
    ,@('E18', '0.66%')
    ,@('D19', '2.188')
    ,@('E19', '0.26%')
    ,@('E21', '0.47%')
    ,@('D22', '4.089')
    ,@('E22', '0.03%')
    ,@('D23', '0.04254')
    ,@('E23', '0.03%')
    ,@('D24', '0.001181')
    ,@('E24', '-3.45%')
    ,@('D25', '0.004052')
    ,@('E25', '6.59%')
    ,@('E26', '0.03%')
    ,@('E27', '4.10%')
    ,@('D40', '0.03975')
    ,@('E40', '2.83%')
    ,@('D41', '0.1113')
    ,@('E41', '0.03%')
    ,@('D42', '0.004140')
    ,@('E42', '1.46%')
    ,@('D43', '0.002142')
    ,@('E43', '-3.36%')
    ,@('E44', '-18.80%')
    ,@('D45', '0.00005164')
    ,@('E45', '0.37%')
    ,@('E46', '0.03%')
    ,@('D48', '0.2582')
    ,@('E48', '90.52%')
    ,@('D49', '0.00002101')
    ,@('E49', '0.03%')
    ,@('D50', '0.0002001')
    ,@('E50', '0.03%')
)

foreach ($u in $updates) {
    $cell = $u[0]
    $newValue = $u[1]
    $rng = $ws.Range($cell)
    # Force text interpretation so numeric-looking strings (e.g. "261.53")
    # and percentages (e.g. "0.24%") are stored as literal text, matching
    # the original inline-string cells instead of being parsed as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    # Reset the cell style back to the workbook default so no stray
    # "Text" number-format style is left attached to the cell.
    $rng.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
